$d = $word.ActiveDocument

# Step 1: remove the paragraph-mark-only rFonts hint="eastAsia" from paragraph 1
# by deleting paragraph 1 (text + mark) and re-inserting an identical paragraph
# whose pPr/rPr lacks the rFonts element.
$p1 = $d.Paragraphs(1)
$p1Range = $p1.Range
$p1Range.Delete()

$newPara1Xml = @'
<?xml version='1.0'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p w14:paraId="010880B8" w14:textId="667C570E" w:rsidR="00560FF7" w:rsidRPr="00560FF7" w:rsidRDefault="00560FF7" w:rsidP="00560FF7"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Файл </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">features -&gt; </w:t></w:r><w:r w:rsidRPr="00560FF7"><w:t>features_NHANES3_HDTrain_all</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> - &gt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>видалит</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>и</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>‘</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US"/></w:rPr><w:t>age</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>’</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$insPoint1 = $d.Range(0, 0)
$insPoint1.InsertXML($newPara1Xml)

# Step 2: insert the four new paragraphs (page break, heading, code, blank)
# right after paragraph 1 and before the pre-existing trailing empty paragraph.
$afterPara1 = $d.Paragraphs(1).Range.End
$insPoint2 = $d.Range($afterPara1, $afterPara1)

$newParasXml = @'
<?xml version='1.0'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>def</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>biological_age_of_each_cluster</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t xml:space="preserve"># </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>summ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> += </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ages_train_dataframe</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>['</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Age</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>'].</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>values</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>person_index</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>]</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$insPoint2.InsertXML($newParasXml)
